# Commit: "Remove mention of standard form certificate from terms"
#
# The published-URL sentence currently reads:
#   "...published at https://{{{domain}}}, substituting this page for the
#   "standard form certificate" those terms refer to."
# and must become:
#   "...published at https://{{{domain}}}."
#
# The three trailing runs (", ", the "substituting..." sentence, and the
# final ".") all share the same InternetLink character style, so a single
# Find/Replace across that stretch of text collapses them back down to one
# run, which is exactly what the target markup shows. The search
# deliberately starts right after "{{{domain}}}" (rather than including it)
# so the bookmark that ends immediately before that run is left anchored in
# its original place instead of being dragged along by the replace.

$d = $word.ActiveDocument

$openQuote  = [char]0x201C
$closeQuote = [char]0x201D

$oldText = ", substituting this page for the " + $openQuote + `
    "standard form certificate" + $closeQuote + " those terms refer to."
$newText = "."

$found = $d.Content.Find.Execute(
    $oldText, $true, $false, $false, $false, $false,
    $true, 1, $false, $newText, 2)

if (-not $found) {
    throw "Could not find the standard-form-certificate sentence to replace."
}

# --- Best-effort style-table touch-ups that also rode along in the diff ---
# The authoring tool's re-save minted a new (empty) character style,
# "ListLabel 6", alongside the pre-existing "ListLabel 1..5" styles. Add the
# closest equivalent via the Styles collection (character style, quick-style
# flagged) even though the exact byte-level placement/markup of a freshly
# minted style is controlled by the host, not script-visible.
$styles = $d.Styles
$hasListLabel6 = $false
foreach ($s in $styles) {
    if ($s.NameLocal -eq "ListLabel 6") { $hasListLabel6 = $true }
}
if (-not $hasListLabel6) {
    $listLabel6 = $styles.Add("ListLabel 6", 2)
    $listLabel6.QuickStyle = $true
}
